$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.167731165885925
$ws.Range("B1").Value = 2.43769383430481
$ws.Range("D1").Value = 2.368287563323975
$ws.Range("E1").Value = 1.234218716621399
